$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Complete" column to the table / worksheet -----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null

# Header
$ws.Range("F1").Value = "Complete"

# --- Fix up / split existing quest-complete dialogue text --------------------
$ws.Range("D2").Value = "믿음직한 당신! 혹시 내 볼펜은 못 봤어..? 내가 아끼는 건데…"
$ws.Range("E2").Value = "내 볼펜이다! 역시 찾아줄 줄 알았어."
$ws.Range("E3").Value = "암호를 풀어내다니… 영리한 인간이구만. "

# --- New "Complete" column dialogue values -----------------------------------
$ws.Range("F2").Value = "퇴근시간을 놓쳐버렸어.. 밤엔 그들이 나오는데…"
$ws.Range("F3").Value = "그래! 그 메시지를 잘 기억하게나..!"
$ws.Range("F4").Value = "악마 사원들 보면 인간 사원들과 별로 다르지 않아요. 모두 힘들어보이거든요."
$ws.Range("F5").Value = "오…오늘..저..점심을..먹고…난..뒤로…모…몸이..이상해…"
$ws.Range("D6").Value = "우리는 시키는 일을 할 뿐이야 히히"
$ws.Range("F7").Value = "캬~ 얼른 퇴근하고 소설을 읽고싶구만!"

# --- Selection cosmetics (matches the saved file's last selection) ----------
$ws.Range("F7").Select() | Out-Null
